$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new columns at their final target positions, left-to-right,
# so each insertion index is valid against the already-settled state.
$ws.Columns("B").Insert()
$ws.Columns("E").Insert()
$ws.Columns("F").Insert()
$ws.Columns("G").Insert()
$ws.Columns("I").Insert()
$ws.Columns("J").Insert()
$ws.Columns("K").Insert()
$ws.Columns("L").Insert()
$ws.Columns("N").Insert()
$ws.Columns("P").Insert()
$ws.Columns("R").Insert()
$ws.Columns("T").Insert()
$ws.Columns("U").Insert()
$ws.Columns("W").Insert()

# Header row (row 1): set names for the newly inserted columns.
$ws.Range("B1").Value = "Alexis Rainey"
$ws.Range("E1").Value = "Curley"
$ws.Range("F1").Value = "Doyle"
$ws.Range("G1").Value = "Espona"
$ws.Range("I1").Value = "Hackman"
$ws.Range("J1").Value = "Holzman"
$ws.Range("K1").Value = "Hughes"
$ws.Range("L1").Value = "Johnson"
$ws.Range("N1").Value = "McFadden"
$ws.Range("P1").Value = "Myers"
$ws.Range("R1").Value = "Reilly"
$ws.Range("T1").Value = "Streib"
$ws.Range("U1").Value = "Tollaksen"
$ws.Range("W1").Value = "Yanovich"

# New data columns inherited column-A formatting on insert; clear it
# for the data rows (2-6) so only the header row keeps the bordered style.
$ws.Range("B2:B6").ClearFormats()
$ws.Range("E2:E6").ClearFormats()
$ws.Range("F2:F6").ClearFormats()
$ws.Range("G2:G6").ClearFormats()
$ws.Range("I2:I6").ClearFormats()
$ws.Range("J2:J6").ClearFormats()
$ws.Range("K2:K6").ClearFormats()
$ws.Range("L2:L6").ClearFormats()
$ws.Range("N2:N6").ClearFormats()
$ws.Range("P2:P6").ClearFormats()
$ws.Range("R2:R6").ClearFormats()
$ws.Range("T2:T6").ClearFormats()
$ws.Range("U2:U6").ClearFormats()
$ws.Range("W2:W6").ClearFormats()

# Fill in the data values for rows 2-6 across all columns B:W.
$ws.Range("B2").Value = 3.7702
$ws.Range("C2").Value = 3.5017
$ws.Range("D2").Value = 3.7225
$ws.Range("F2").Value = 4.2893
$ws.Range("G2").Value = 4.4326
$ws.Range("H2").Value = 3.3586
$ws.Range("I2").Value = 3.6899
$ws.Range("J2").Value = 3.6008
$ws.Range("L2").Value = 3.9427
$ws.Range("M2").Value = 4.1038
$ws.Range("N2").Value = 4.2551
$ws.Range("O2").Value = 3.6064
$ws.Range("P2").Value = 3.6876
$ws.Range("Q2").Value = 3.5987
$ws.Range("S2").Value = 3.6262
$ws.Range("U2").Value = 4.2828
$ws.Range("V2").Value = 3.7954
$ws.Range("W2").Value = 3.3294
$ws.Range("B3").Value = 2.8984
$ws.Range("C3").Value = 2.3147
$ws.Range("D3").Value = 1.6977
$ws.Range("E3").Value = 1.7114
$ws.Range("F3").Value = 2.577
$ws.Range("G3").Value = 2.3936
$ws.Range("H3").Value = 1.559
$ws.Range("I3").Value = 2.1728
$ws.Range("J3").Value = 2.3316
$ws.Range("K3").Value = 1.6503
$ws.Range("L3").Value = 3.3331
$ws.Range("M3").Value = 2.8072
$ws.Range("N3").Value = 2.3012
$ws.Range("O3").Value = 2.7648
$ws.Range("P3").Value = 1.9123
$ws.Range("Q3").Value = 2.0339
$ws.Range("R3").Value = 2.6011
$ws.Range("S3").Value = 2.7599
$ws.Range("U3").Value = 2.4253
$ws.Range("V3").Value = 2.0448
$ws.Range("W3").Value = 2.2357
$ws.Range("B4").Value = 3.1076
$ws.Range("C4").Value = 2.5512
$ws.Range("D4").Value = 2.6238
$ws.Range("E4").Value = 1.447
$ws.Range("F4").Value = 3.6183
$ws.Range("G4").Value = 2.742
$ws.Range("H4").Value = 3.2772
$ws.Range("I4").Value = 3.2636
$ws.Range("J4").Value = 3.1381
$ws.Range("K4").Value = 2.474
$ws.Range("L4").Value = 3.8067
$ws.Range("M4").Value = 3.6149
$ws.Range("N4").Value = 3.0526
$ws.Range("O4").Value = 3.913
$ws.Range("P4").Value = 3.0377
$ws.Range("Q4").Value = 2.8756
$ws.Range("R4").Value = 3.1682
$ws.Range("S4").Value = 3.1551
$ws.Range("U4").Value = 3.7154
$ws.Range("V4").Value = 3.0091
$ws.Range("B5").Value = 2.857
$ws.Range("C5").Value = 2.5587
$ws.Range("E5").Value = 1.5596
$ws.Range("F5").Value = 3.1613
$ws.Range("G5").Value = 3.2181
$ws.Range("H5").Value = 2.1635
$ws.Range("I5").Value = 2.3426
$ws.Range("J5").Value = 2.3795
$ws.Range("K5").Value = 2.1696
$ws.Range("L5").Value = 1.8492
$ws.Range("M5").Value = 2.9072
$ws.Range("N5").Value = 2.9156
$ws.Range("O5").Value = 2.7582
$ws.Range("Q5").Value = 2.6653
$ws.Range("R5").Value = 2.4971
$ws.Range("S5").Value = 2.7622
$ws.Range("T5").Value = 1.2709
$ws.Range("U5").Value = 3.0247
$ws.Range("V5").Value = 2.6411
$ws.Range("B6").Value = 2.7425
$ws.Range("C6").Value = 2.1367
$ws.Range("E6").Value = 1.0674
$ws.Range("F6").Value = 2.5891
$ws.Range("G6").Value = 2.7828
$ws.Range("H6").Value = 2.1273
$ws.Range("I6").Value = 2.5152
$ws.Range("J6").Value = 1.5082
$ws.Range("K6").Value = 0.8742
$ws.Range("L6").Value = 2.6764
$ws.Range("M6").Value = 1.3475
$ws.Range("N6").Value = 1.8795
$ws.Range("O6").Value = 2.7213
$ws.Range("Q6").Value = 1.9802
$ws.Range("R6").Value = 2.1594
$ws.Range("S6").Value = 1.9253
$ws.Range("T6").Value = 1.0585
$ws.Range("V6").Value = 1.8846
